# Add the "2022-Q4" sheet before the current "2022-Q3" sheet (position 2),
# and insert a new top row into the "总计" (total) summary sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert a new worksheet named "2022-Q4" right before "2022-Q3".
# ---------------------------------------------------------------------
$q3Before = $wb.Worksheets.Item(2)
$newSheet = $wb.Worksheets.Add($q3Before)
$newSheet.Name = "2022-Q4"

# Re-fetch the "2022-Q3" worksheet by its (now shifted) position so we get a
# live reference to the sheet that still holds the data/styles we want to
# clone (this runtime tracks worksheet variables by position, not identity).
$q3 = $wb.Worksheets.Item(3)

# Clone header row + first 15 rows worth of formatting (2022-Q4 only needs
# 14 data rows) from "2022-Q3" into the new sheet. This brings over the
# bold/bordered/centered style used for the header row and column A index.
$q3.Range("A1:H15").Copy($newSheet.Range("A1"))

# A "plain" (unstyled) cell in the freshly pasted block, used below to reset
# the style of cells whose value we re-enter with a leading quote (to force
# text instead of number) back to the default/unstyled look.
$plainStyle = $newSheet.Cells.Item(2, 3).Style

function Set-TextValue($cell, $text) {
    $cell.Value = "'" + $text
    $cell.Style = $plainStyle
}

# ---------------------------------------------------------------------
# 2) Populate the 2022-Q4 fund holdings data (header row is already correct
#    after the copy above, so only data rows 2-15 are written).
# ---------------------------------------------------------------------
$data = @(
    @("009556", "兴全合丰三年持有期混合", "65.69", "94.32", "3.96", "2.6013", 10),
    @("900090", "中信卓越成长两年持有期混合B", "56.19", "93.11", "3.23", "1.8149", 9),
    @("005644", "广发沪港深行业龙头混合", "14.31", "89.49", "4.58", "0.6554", 7),
    @("900010", "中信卓越成长两年持有期混合A", "16.47", "93.11", "3.23", "0.5320", 9),
    @("005228", "汇添富港股通专注成长混合", "7.51", "85.64", "5.02", "0.3770", 4),
    @("013123", "汇添富精选核心优势一年持有混合A", "5.58", "83.43", "4.66", "0.2600", 6),
    @("900100", "中信卓越成长两年持有期混合C", "4.96", "93.11", "3.23", "0.1602", 9),
    @("006696", "汇添富研究优选灵活配置混合", "3.37", "76.35", "2.98", "0.1004", 8),
    @("013550", "汇添富品牌价值一年持有混合A", "2.24", "75.70", "4.31", "0.0965", 5),
    @("010480", "汇添富稳进双盈一年持有期混合", "7.92", "23.94", "0.65", "0.0515", 9),
    @("013367", "汇添富多元价值发现混合A", "0.77", "65.40", "3.02", "0.0233", 7),
    @("013124", "汇添富精选核心优势一年持有混合C", "0.27", "83.43", "4.66", "0.0126", 6),
    @("013551", "汇添富品牌价值一年持有混合C", "0.28", "75.70", "4.31", "0.0121", 5),
    @("013368", "汇添富多元价值发现混合C", "0.27", "65.40", "3.02", "0.0082", 7)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]

    # Column B: fund code -> always digits, must stay text.
    Set-TextValue $newSheet.Cells.Item($r, 2) $row[0]

    # Column C: fund name -> never numeric, safe to assign directly.
    $newSheet.Cells.Item($r, 3).Value = $row[1]

    # Columns D-G: numeric-looking decimal text, must stay text.
    Set-TextValue $newSheet.Cells.Item($r, 4) $row[2]
    Set-TextValue $newSheet.Cells.Item($r, 5) $row[3]
    Set-TextValue $newSheet.Cells.Item($r, 6) $row[4]
    Set-TextValue $newSheet.Cells.Item($r, 7) $row[5]

    # Column H: rank -> real number.
    $newSheet.Cells.Item($r, 8).Value = $row[6]
}

# ---------------------------------------------------------------------
# 3) Update the "总计" summary sheet: insert a new row for 2022-Q4 above the
#    existing 2022-Q3 row, pushing the rest down.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()

$idxCell = $summary.Cells.Item(2, 1)
$idxCell.Value = 0
$idxCell.Style = $summary.Cells.Item(3, 1).Style

$summary.Cells.Item(2, 2).Value = "2022-Q4"
$summary.Cells.Item(2, 3).Value = 14
$summary.Cells.Item(2, 4).Value = 6.71
